$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1800.2858
$ws.Range("J40").Value = 1560.4
$ws.Range("L40").Value = 1560.4
$ws.Range("N40").Value = -1910.4

$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376

$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880

$ws.Range("H129").Value = 1054.4468
$ws.Range("J129").Value = 1180.8
$ws.Range("L129").Value = 3542.4
$ws.Range("N129").Value = -13542.4

$ws.Range("H132").Value = 10955.454
$ws.Range("I132").Value = 10955.454
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 32866.362
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 1404.5
$ws.Range("I137").Value = 1265.5
$ws.Range("J137").Value = 1456.625
$ws.Range("K137").Value = 3796.5
$ws.Range("L137").Value = 4369.875
$ws.Range("M137").Value = -1246.5
$ws.Range("N137").Value = -9469.875

$ws.Range("H141").Value = 3379.2258
$ws.Range("I141").Value = 1877.125
$ws.Range("J141").Value = 8529.286
$ws.Range("K141").Value = 5631.375
$ws.Range("L141").Value = 25587.858
$ws.Range("M141").Value = -451.375
$ws.Range("N141").Value = -35947.858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1108.1428
$ws.Range("I2").Value = 1202.75
$ws.Range("K2").Value = 1202.75
$ws.Range("M2").Value = -1089.75

$ws.Range("H32").Value = 459044.03
$ws.Range("I32").Value = 666940.1
$ws.Range("J32").Value = 12067.35
$ws.Range("K32").Value = 666940.1
$ws.Range("L32").Value = 12067.35
$ws.Range("M32").Value = -666653.1
$ws.Range("N32").Value = -12641.35

$ws.Range("H116").Value = 1108.1428
$ws.Range("I116").Value = 1202.75
$ws.Range("K116").Value = 1202.75
$ws.Range("M116").Value = 1091.25

$ws.Range("H122").Value = 1529.9375
$ws.Range("I122").Value = 1379.9
$ws.Range("J122").Value = 1780
$ws.Range("K122").Value = 4139.700000000001
$ws.Range("L122").Value = 5340
$ws.Range("M122").Value = -1689.700000000001
$ws.Range("N122").Value = -10240

$ws.Range("H132").Value = 6610.148
$ws.Range("I132").Value = 6351.8945
$ws.Range("J132").Value = 7223.5
$ws.Range("K132").Value = 19055.6835
$ws.Range("L132").Value = 21670.5
$ws.Range("M132").Value = -16525.6835
$ws.Range("N132").Value = -26730.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1108.1428
$ws.Range("I3").Value = 1202.75
$ws.Range("K3").Value = 1202.75
$ws.Range("M3").Value = -1088.75

$ws.Range("H64").Value = 998.25
$ws.Range("J64").Value = 1043.4
$ws.Range("L64").Value = 1043.4
$ws.Range("N64").Value = -1493.4

$ws.Range("H67").Value = 998.25
$ws.Range("J67").Value = 1043.4
$ws.Range("L67").Value = 1043.4
$ws.Range("N67").Value = -2603.4

$ws.Range("H99").Value = 849.2
$ws.Range("I99").Value = 743.6842
$ws.Range("J99").Value = 1183.3334
$ws.Range("K99").Value = 743.6842
$ws.Range("L99").Value = 1183.3334
$ws.Range("M99").Value = 754.3158
$ws.Range("N99").Value = -4179.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1627.619
$ws.Range("I31").Value = 1267.3103
$ws.Range("J31").Value = 2431.3845
$ws.Range("K31").Value = 1267.3103
$ws.Range("L31").Value = 2431.3845
$ws.Range("M31").Value = -972.3103000000001
$ws.Range("N31").Value = -3021.3845

$ws.Range("H34").Value = 1627.619
$ws.Range("I34").Value = 1267.3103
$ws.Range("J34").Value = 2431.3845
$ws.Range("K34").Value = 1267.3103
$ws.Range("L34").Value = 2431.3845
$ws.Range("M34").Value = -1065.3103
$ws.Range("N34").Value = -2835.3845

$ws.Range("H50").Value = 13197.2
$ws.Range("J50").Value = 13197.2
$ws.Range("L50").Value = 13197.2
$ws.Range("N50").Value = -14447.2

$ws.Range("H122").Value = 2235.2666
$ws.Range("I122").Value = 2435.75
$ws.Range("J122").Value = 1433.3334
$ws.Range("K122").Value = 7307.25
$ws.Range("L122").Value = 4300.0002
$ws.Range("M122").Value = -4857.25
$ws.Range("N122").Value = -9200.0002

$ws.Range("H132").Value = 9806246
$ws.Range("I132").Value = 1667.1666
$ws.Range("K132").Value = 5001.4998
$ws.Range("M132").Value = -2471.4998

$ws.Range("H134").Value = 1173.8182
$ws.Range("I134").Value = 768
$ws.Range("K134").Value = 2304
$ws.Range("M134").Value = 231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 874.7857
$ws.Range("I68").Value = 856.2857
$ws.Range("J68").Value = 876.2088
$ws.Range("K68").Value = 2568.8571
$ws.Range("L68").Value = 2628.6264
$ws.Range("M68").Value = -1757.8571
$ws.Range("N68").Value = -4250.6264

$ws.Range("H71").Value = 874.7857
$ws.Range("I71").Value = 856.2857
$ws.Range("J71").Value = 876.2088
$ws.Range("K71").Value = 7706.571300000001
$ws.Range("L71").Value = 7885.8792
$ws.Range("M71").Value = -3650.571300000001
$ws.Range("N71").Value = -15997.8792

$ws.Range("H107").Value = 1439.0408
$ws.Range("J107").Value = 2999.8096
$ws.Range("L107").Value = 8999.4288
$ws.Range("N107").Value = -12839.4288

$ws.Range("H113").Value = 691.5143
$ws.Range("I113").Value = 396.85184
$ws.Range("J113").Value = 1686
$ws.Range("K113").Value = 1190.55552
$ws.Range("L113").Value = 5058
$ws.Range("M113").Value = 979.4444800000001
$ws.Range("N113").Value = -9398

$ws.Range("H131").Value = 1043.7727
$ws.Range("I131").Value = 840.9091
$ws.Range("J131").Value = 1111.3939
$ws.Range("K131").Value = 2522.7273
$ws.Range("L131").Value = 3334.1817
$ws.Range("M131").Value = 2517.2727
$ws.Range("N131").Value = -13414.1817

$ws.Range("H132").Value = 3571.5
$ws.Range("I132").Value = 2647.087
$ws.Range("J132").Value = 4178.971
$ws.Range("K132").Value = 23823.783
$ws.Range("L132").Value = 37610.73899999999
$ws.Range("M132").Value = -21293.783
$ws.Range("N132").Value = -42670.73899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2104.375
$ws.Range("I113").Value = 1830.5
$ws.Range("J113").Value = 2378.25
$ws.Range("K113").Value = 1830.5
$ws.Range("L113").Value = 2378.25
$ws.Range("M113").Value = 339.5
$ws.Range("N113").Value = -6718.25

$ws.Range("H132").Value = 3156.7
$ws.Range("I132").Value = 2821
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 8463
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -5933
$ws.Range("N132").Value = -18558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5125
$ws.Range("I7").Value = 5166.6665
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 5166.6665
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -5054.6665
$ws.Range("N7").Value = -5224

$ws.Range("H126").Value = 5125
$ws.Range("I126").Value = 5166.6665
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 15499.9995
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -13029.9995
$ws.Range("N126").Value = -19940

$ws.Range("H132").Value = 3356.2
$ws.Range("I132").Value = 3004.2903
$ws.Range("J132").Value = 3930.3684
$ws.Range("K132").Value = 9012.8709
$ws.Range("L132").Value = 11791.1052
$ws.Range("M132").Value = -6482.8709
$ws.Range("N132").Value = -16851.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1482.7858
$ws.Range("I126").Value = 1260.4445
$ws.Range("J126").Value = 1883
$ws.Range("K126").Value = 3781.3335
$ws.Range("L126").Value = 5649
$ws.Range("M126").Value = -1311.3335
$ws.Range("N126").Value = -10589

$ws.Range("H136").Value = 2628.4565
$ws.Range("I136").Value = 2528.36
$ws.Range("J136").Value = 2747.6191
$ws.Range("K136").Value = 7585.08
$ws.Range("L136").Value = 8242.8573
$ws.Range("M136").Value = -5035.08
$ws.Range("N136").Value = -13342.8573
